# Update movie name values: the two "Avatar" variants are consolidated
# into a single, consistently-spelled entry "Avatar: way of water".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dim = $ws.Range("A1:F46")
$rowCount = $dim.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 3)   # column C = movieName
    $val = $cell.Value2
    if ($val -eq "Avatar, way of water" -or $val -eq "Avatar: the way of water") {
        $cell.Value = "Avatar: way of water"
    }
}

# Reflect the author's final selection/view state in the sheet.
$ws.Range("C5").Select()
